# Generate Report for Handoff
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps for the 629fd7b1... and 74b534fd... handoff rows.
# - Mark those same rows (and the other "Ready for handoff" rows) with the
#   "ht" priority on both the zh-cn and de-de localization sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows in each status table that hold "Ready for handoff" files.
$rows = @(8, 9, 11, 12, 13, 14)

foreach ($r in $rows) {
    # Overview!G -> "Latest HO Xliff Generate Date" for de-de rolls the
    # handoff timestamp forward from 04:20:25 to 04:20:41.
    $overview.Range("G$r").Value = "2016-08-30 04:20:41"

    # de-de!H -> "Latest Handoff Datetime" mirrors the Overview column.
    $dede.Range("H$r").Value = "2016-08-30 04:20:41"

    # zh-cn!H -> "Latest Handoff Datetime" rolls forward from 04:20:20 to 04:20:36.
    $zhcn.Range("H$r").Value = "2016-08-30 04:20:36"

    # Priority column E on both localization sheets now reports "ht".
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
